$wb = $excel.ActiveWorkbook

# ALC row 62 (hunk 0)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 78620.266
$ws.Range("I62").Value = 103409.45
$ws.Range("J62").Value = 10450
$ws.Range("K62").Value = 103409.45
$ws.Range("L62").Value = 10450
$ws.Range("M62").Value = -102785.45
$ws.Range("N62").Value = -11698

# ALC row 65 (hunk 1)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 78620.266
$ws.Range("I65").Value = 103409.45
$ws.Range("J65").Value = 10450
$ws.Range("K65").Value = 517047.25
$ws.Range("L65").Value = 52250
$ws.Range("M65").Value = -513927.25
$ws.Range("N65").Value = -58490

# ALC row 116 (hunk 2)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 5000.75
$ws.Range("I116").Value = 6175
$ws.Range("K116").Value = 6175
$ws.Range("M116").Value = -2733

# ARM row 2 (hunk 3)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2277.6086
$ws.Range("I2").Value = 2298.25
$ws.Range("J2").Value = 2230.4285
$ws.Range("K2").Value = 2298.25
$ws.Range("L2").Value = 2230.4285
$ws.Range("M2").Value = -2185.25
$ws.Range("N2").Value = -2456.4285

# ARM row 4 (hunk 4)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 282.1111
$ws.Range("I4").Value = 267.8
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 267.8
$ws.Range("L4").Value = 300
$ws.Range("M4").Value = -151.8
$ws.Range("N4").Value = -532

# ARM row 88 (hunk 5)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 4500
$ws.Range("I88").Value = 4000
$ws.Range("J88").Value = 5000
$ws.Range("K88").Value = 4000
$ws.Range("L88").Value = 5000
$ws.Range("M88").Value = -3594
$ws.Range("N88").Value = -5812

# ARM row 91 (hunk 6)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 4500
$ws.Range("I91").Value = 4000
$ws.Range("J91").Value = 5000
$ws.Range("K91").Value = 4000
$ws.Range("L91").Value = 5000
$ws.Range("M91").Value = -2596
$ws.Range("N91").Value = -7808

# ARM row 116 (hunk 7)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2277.6086
$ws.Range("I116").Value = 2298.25
$ws.Range("J116").Value = 2230.4285
$ws.Range("K116").Value = 2298.25
$ws.Range("L116").Value = 2230.4285
$ws.Range("M116").Value = -4.25
$ws.Range("N116").Value = -6818.4285

# BSM row 3 (hunk 8)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2277.6086
$ws.Range("I3").Value = 2298.25
$ws.Range("J3").Value = 2230.4285
$ws.Range("K3").Value = 2298.25
$ws.Range("L3").Value = 2230.4285
$ws.Range("M3").Value = -2184.25
$ws.Range("N3").Value = -2458.4285

# BSM row 5 (hunk 9)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 56
$ws.Range("I5").Value = 56
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 56
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 57
$ws.Range("N5").ClearContents()

# BSM row 63 (hunk 10)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

# BSM row 66 (hunk 11)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

# BSM row 69 (hunk 12)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H69").Value = 34333.332
$ws.Range("J69").Value = 34333.332
$ws.Range("L69").Value = 34333.332
$ws.Range("N69").Value = -35955.332

# BSM row 72 (hunk 13)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H72").Value = 34333.332
$ws.Range("J72").Value = 34333.332
$ws.Range("L72").Value = 102999.996
$ws.Range("N72").Value = -111111.996

# BSM row 75 (hunk 14)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 30499.75
$ws.Range("I75").Value = 12000
$ws.Range("J75").Value = 36666.332
$ws.Range("K75").Value = 12000
$ws.Range("L75").Value = 36666.332
$ws.Range("M75").Value = -11064
$ws.Range("N75").Value = -38538.332

# BSM row 78 (hunk 15)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H78").Value = 30499.75
$ws.Range("I78").Value = 12000
$ws.Range("J78").Value = 36666.332
$ws.Range("K78").Value = 36000
$ws.Range("L78").Value = 109998.996
$ws.Range("M78").Value = -31320
$ws.Range("N78").Value = -119358.996

# BSM row 94 (hunk 16)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1176.5
$ws.Range("I94").Value = 859.7143
$ws.Range("K94").Value = 859.7143
$ws.Range("M94").Value = -408.7143

# BSM row 134 (hunk 17)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 77465.14999999999
$ws.Range("I134").Value = 3564.375
$ws.Range("K134").Value = 10693.125
$ws.Range("M134").Value = -8158.125

# CRP row 20 (hunk 18)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 47151.6
$ws.Range("J20").Value = 47151.6
$ws.Range("L20").Value = 47151.6
$ws.Range("N20").Value = -47623.6

# CRP row 30 (hunk 19)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H30").Value = 47151.6
$ws.Range("J30").Value = 47151.6
$ws.Range("L30").Value = 47151.6
$ws.Range("N30").Value = -47333.6

# CRP row 128 (hunk 20)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H128").Value = 47151.6
$ws.Range("J128").Value = 47151.6
$ws.Range("L128").Value = 47151.6
$ws.Range("N128").Value = -57111.6

# CUL row 34 (hunk 21)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1409
$ws.Range("I34").Value = 135
$ws.Range("J34").Value = 1727.5
$ws.Range("K34").Value = 405
$ws.Range("L34").Value = 5182.5
$ws.Range("M34").Value = -321
$ws.Range("N34").Value = -5350.5

# CUL row 69 (hunk 22)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 951.2
$ws.Range("I69").Value = 744.5714
$ws.Range("J69").Value = 1433.3334
$ws.Range("K69").Value = 2233.7142
$ws.Range("L69").Value = 4300.0002
$ws.Range("M69").Value = -1422.7142
$ws.Range("N69").Value = -5922.0002

# CUL row 72 (hunk 23)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H72").Value = 951.2
$ws.Range("I72").Value = 744.5714
$ws.Range("J72").Value = 1433.3334
$ws.Range("K72").Value = 6701.1426
$ws.Range("L72").Value = 12900.0006
$ws.Range("M72").Value = -2645.1426
$ws.Range("N72").Value = -21012.0006

# CUL row 87 (hunk 24)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 13899.8
$ws.Range("I87").Value = 5571.4287
$ws.Range("K87").Value = 16714.2861
$ws.Range("M87").Value = -15466.2861

# CUL row 90 (hunk 25)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 13899.8
$ws.Range("I90").Value = 5571.4287
$ws.Range("K90").Value = 50142.85830000001
$ws.Range("M90").Value = -43902.85830000001

# CUL row 121 (hunk 26)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 409.66666
$ws.Range("J121").Value = 499.5
$ws.Range("L121").Value = 1498.5
$ws.Range("N121").Value = -4118.5

# GSM row 102 (hunk 27)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1315.04
$ws.Range("I102").Value = 1075.6666
$ws.Range("J102").Value = 1930.5714
$ws.Range("K102").Value = 1075.6666
$ws.Range("L102").Value = 1930.5714
$ws.Range("M102").Value = 546.3334
$ws.Range("N102").Value = -5174.5714

# GSM row 122 (hunk 28)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2060.025
$ws.Range("I122").Value = 1591.92
$ws.Range("J122").Value = 2840.2
$ws.Range("K122").Value = 4775.76
$ws.Range("L122").Value = 8520.599999999999
$ws.Range("M122").Value = -2325.76
$ws.Range("N122").Value = -13420.6

# GSM row 123 (hunk 29)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 25702.285
$ws.Range("J123").Value = 25702.285
$ws.Range("L123").Value = 25702.285
$ws.Range("N123").Value = -30602.285

# GSM row 132 (hunk 30)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3123.125
$ws.Range("I132").Value = 2597.2
$ws.Range("K132").Value = 7791.599999999999
$ws.Range("M132").Value = -5261.599999999999

# LTW row 7 (hunk 31)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2100.15
$ws.Range("I7").Value = 2021.6428
$ws.Range("J7").Value = 2283.3333
$ws.Range("K7").Value = 2021.6428
$ws.Range("L7").Value = 2283.3333
$ws.Range("M7").Value = -1909.6428
$ws.Range("N7").Value = -2507.3333

# LTW row 40 (hunk 32)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2080.261
$ws.Range("I40").Value = 2068.2104
$ws.Range("J40").Value = 2137.5
$ws.Range("K40").Value = 2068.2104
$ws.Range("L40").Value = 2137.5
$ws.Range("M40").Value = -1932.2104
$ws.Range("N40").Value = -2409.5

# LTW row 61 (hunk 33)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6309.2085
$ws.Range("I61").Value = 6496.2383
$ws.Range("K61").Value = 6496.2383
$ws.Range("M61").Value = -6294.2383

# LTW row 68 (hunk 34)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3533.3333
$ws.Range("I68").Value = 3457.1428
$ws.Range("K68").Value = 3457.1428
$ws.Range("M68").Value = -2708.1428

# LTW row 71 (hunk 35)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3533.3333
$ws.Range("I71").Value = 3457.1428
$ws.Range("K71").Value = 17285.714
$ws.Range("M71").Value = -13541.714

# LTW row 113 (hunk 36)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 6309.2085
$ws.Range("I113").Value = 6496.2383
$ws.Range("K113").Value = 6496.2383
$ws.Range("M113").Value = -4326.2383

# LTW row 122 (hunk 37)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5115.7334
$ws.Range("I122").Value = 8076.552
$ws.Range("J122").Value = 2345.9355
$ws.Range("K122").Value = 24229.656
$ws.Range("L122").Value = 7037.806500000001
$ws.Range("M122").Value = -21779.656
$ws.Range("N122").Value = -11937.8065

# LTW row 126 (hunk 38)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2100.15
$ws.Range("I126").Value = 2021.6428
$ws.Range("J126").Value = 2283.3333
$ws.Range("K126").Value = 6064.928400000001
$ws.Range("L126").Value = 6849.999899999999
$ws.Range("M126").Value = -3594.928400000001
$ws.Range("N126").Value = -11789.9999

# LTW row 136 (hunk 39)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2470.9285
$ws.Range("I136").Value = 1122.5385
$ws.Range("K136").Value = 3367.6155
$ws.Range("M136").Value = -817.6155000000003
